$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 1: turn the stray duplicate data row into a proper header row ---
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").Value = "name"

$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").Value = "capacity"

$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("D1").Value = "owner"

$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4160
$ws.Range("E1").Borders.LineStyle = 1
$ws.Range("E1").Value = "register_date"

$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1
$ws.Range("F1").Value = "register_reason"

$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").VerticalAlignment = -4160
$ws.Range("G1").Borders.LineStyle = 1
$ws.Range("G1").Value = "acquire_value"

$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1
$ws.Range("H1").Value = "property_category"

$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160
$ws.Range("I1").Borders.LineStyle = 1
$ws.Range("I1").Value = "category"

$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1
$ws.Range("J1").Value = "date"

$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("K1").VerticalAlignment = -4160
$ws.Range("K1").Borders.LineStyle = 1
$ws.Range("K1").Value = "legislator_name"

$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("L1").VerticalAlignment = -4160
$ws.Range("L1").Borders.LineStyle = 1
$ws.Range("L1").Value = "legislator_id"

$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").HorizontalAlignment = -4108
$ws.Range("M1").VerticalAlignment = -4160
$ws.Range("M1").Borders.LineStyle = 1
$ws.Range("M1").Value = "source_file"

$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").HorizontalAlignment = -4108
$ws.Range("N1").VerticalAlignment = -4160
$ws.Range("N1").Borders.LineStyle = 1
$ws.Range("N1").Value = "index"

# --- Row 2 (car #35, 中華): fill in the new property/legislator columns ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2013-12-26"
$ws.Range("K2").Value = "李桐豪"
$ws.Range("L2").Value = 896
$ws.Range("M2").Value = "tmp2e9d1"
$ws.Range("N2").Value = 35

# --- Row 3 (car #36, 福特): fill in the new property/legislator columns ---
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2013-12-26"
$ws.Range("K3").Value = "李桐豪"
$ws.Range("L3").Value = 896
$ws.Range("M3").Value = "tmp2e9d1"
$ws.Range("N3").Value = 36
